$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix mislabeled columns: Nominal -> Ordinal for specific rows
# Row 8 (LandContour), Row 9 (Utilities), Row 63 (GarageQual), Row 64 (GarageCond),
# Row 65 (PavedDrive), Row 72 (PoolQC), Row 73 (Fence)
$ordinalRows = @(8, 9, 63, 64, 65, 72, 73)
foreach ($r in $ordinalRows) {
    $ws.Cells.Item($r, 3).Value = "Ordinal"
}

# Add new row 80 for SalePrice column documentation
$ws.Cells.Item(80, 1).Value = "SalePrice"
$ws.Cells.Item(80, 2).Value = "Sale Price"
$ws.Cells.Item(80, 3).Value = "Discrete"
$ws.Cells.Item(80, 5).Value = "Dollars"

# Give the new "Name" cell (A80) a distinct font style (Arial 10, black) like the header font family
$ws.Cells.Item(80, 1).Font.Name = "Arial"
$ws.Cells.Item(80, 1).Font.Size = 10
$ws.Cells.Item(80, 1).Font.Color = 0

# Update selection/view state to match new active cell
$ws.Cells.Item(80, 3).Select()
